$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: Status moved from "Active" to "Complete"; %Complete moved from 0.4 to 1
$ws.Range("B6").Value = "Complete"
$ws.Range("F6").Value = 1

# New row 12: "EnviroDIY publisher Source" task, Proposed status, with description.
# Write the Description (column G) first so the new shared strings land in the
# same order as the target workbook (description string before task string).
$ws.Range("G12").Value = "Node to convert incoming DIY identifiers to a more descriptive value. (How can we get DIY metadata)."
$ws.Range("A12").Value = "EnviroDIY publisher Source"
$ws.Range("B12").Value = "Proposed"

# Restore the recorded selection state from the authored workbook.
$ws.Range("E17").Select() | Out-Null
